# Automatic update of files.
#
# 1) Column C ("Förändrad") is bumped by one day (46059 -> 46060, i.e.
#    2026-02-06 -> 2026-02-07) for every data row (2..24).
# 2) Rows 10..24 get re-sorted (their A..R payload moves to a different
#    row, carrying its Beteckning-specific hyperlink formulas with it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: bump the "Förändrad" date column for every data row (2..24).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# ---------------------------------------------------------------------
# Step 2: rows 10..24 are re-sorted. Below is the final payload that
# belongs in each row after the sort (everything except column C, which
# was already set to 46060 above for all rows).
# ---------------------------------------------------------------------
$rowsData = @(
    @{ Row=10; A="A 24384-2023"; B=45076; F=$null; G=1; H=0; I=0; J=1; K=0; L=0; M=0; N=0; O=1; P=0; Q=1; R="Kolflarnlav"; HasLinks=$true; HasZ=$false }
    @{ Row=11; A="A 24262-2023"; B=45076; F=$null; G=2.5; H=0; I=1; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=1; R="Tibast"; HasLinks=$true; HasZ=$false }
    @{ Row=12; A="A 24368-2023"; B=45076; F=$null; G=3; H=1; I=1; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=1; R="Plattlummer"; HasLinks=$true; HasZ=$false }
    @{ Row=13; A="A 491-2026"; B=46029.42581018519; F=$null; G=3.2; H=1; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=1; R="Blåsippa"; HasLinks=$true; HasZ=$false }
    @{ Row=14; A="A 49421-2023"; B=45211; F="Kommuner"; G=1.4; H=1; I=0; J=0; K=0; L=1; M=0; N=0; O=1; P=1; Q=1; R="Grönfink"; HasLinks=$true; HasZ=$true }
    @{ Row=15; A="A 24233-2023"; B=45076; F=$null; G=3.7; H=0; I=0; J=1; K=0; L=0; M=0; N=0; O=1; P=0; Q=1; R="Svartvit taggsvamp"; HasLinks=$true; HasZ=$false }
    @{ Row=16; A="A 32508-2022"; B=44782; F="Kommuner"; G=4.7; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=17; A="A 24254-2023"; B=45076; F=$null; G=1.2; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=18; A="A 24257-2023"; B=45076; F=$null; G=1; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=19; A="A 27561-2025"; B=45813.48378472222; F=$null; G=0.7; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=20; A="A 27589-2025"; B=45813.51070601852; F=$null; G=2.7; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=21; A="A 24363-2023"; B=45076; F=$null; G=2.4; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=22; A="A 47173-2022"; B=44852; F=$null; G=1.9; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=23; A="A 35404-2022"; B=44798; F="Kommuner"; G=1.4; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
    @{ Row=24; A="A 31246-2022"; B=44771; F=$null; G=4.3; H=0; I=0; J=0; K=0; L=0; M=0; N=0; O=0; P=0; Q=0; R=$null; HasLinks=$false; HasZ=$false }
)

# Hyperlink columns: column letter -> (folder, file-name suffix, extension)
$linkCols = @(
    @{ Col="S"; Folder="artfynd";          Suffix="artfynd";                Ext="xlsx" },
    @{ Col="T"; Folder="kartor";           Suffix="karta";                  Ext="png"  },
    @{ Col="V"; Folder="klagomål";         Suffix="FSC-klagomål";           Ext="docx" },
    @{ Col="W"; Folder="klagomålsmail";    Suffix="FSC-klagomål mail";      Ext="docx" },
    @{ Col="X"; Folder="tillsyn";          Suffix="tillsynsbegäran";        Ext="docx" },
    @{ Col="Y"; Folder="tillsynsmail";     Suffix="tillsynsbegäran mail";   Ext="docx" }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    $ws.Cells.Item($r, 1).Value = $rd.A          # A - Beteckning
    $ws.Cells.Item($r, 2).Value = $rd.B           # B - Datum

    if ($rd.F -ne $null) {
        $ws.Cells.Item($r, 6).Value = $rd.F       # F - Markägare
    } else {
        $ws.Cells.Item($r, 6).ClearContents()
    }

    $ws.Cells.Item($r, 7).Value  = $rd.G           # G - Area (ha)
    $ws.Cells.Item($r, 8).Value  = $rd.H           # H - Fridlysta
    $ws.Cells.Item($r, 9).Value  = $rd.I           # I - Signalarter
    $ws.Cells.Item($r, 10).Value = $rd.J           # J - NT
    $ws.Cells.Item($r, 11).Value = $rd.K           # K - VU
    $ws.Cells.Item($r, 12).Value = $rd.L           # L - EN
    $ws.Cells.Item($r, 13).Value = $rd.M           # M - CR
    $ws.Cells.Item($r, 14).Value = $rd.N           # N - RE
    $ws.Cells.Item($r, 15).Value = $rd.O           # O - Rödlistade
    $ws.Cells.Item($r, 16).Value = $rd.P           # P - Hotade
    $ws.Cells.Item($r, 17).Value = $rd.Q           # Q - Alla arter

    if ($rd.R -ne $null) {
        $ws.Cells.Item($r, 18).Value = $rd.R       # R - Artnamn
    } else {
        $ws.Cells.Item($r, 18).ClearContents()
    }

    if ($rd.HasLinks) {
        foreach ($lc in $linkCols) {
            $url = 'https://klasma.github.io/Logging_0126/' + $lc.Folder + '/' + $rd.A + ' ' + $lc.Suffix + '.' + $lc.Ext
            $formula = '=HYPERLINK("' + $url + '", "' + $rd.A + '")'
            $ws.Range($lc.Col + $r).Formula = $formula
        }
    }

    if ($rd.HasZ) {
        $url = 'https://klasma.github.io/Logging_0126/fåglar/' + $rd.A + ' prioriterade fågelarter.docx'
        $formula = '=HYPERLINK("' + $url + '", "' + $rd.A + '")'
        $ws.Range("Z" + $r).Formula = $formula
    } else {
        $ws.Range("Z" + $r).ClearContents()
    }
}
